# 7.2.1.1 workbook update
# - Corrects the EN shared string "(in per cent)" -> "(in percent)"
# - Adds a new year column (T) for 2023
# - Revises the 2022 "share" figure and widens the year columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the "(in per cent)" / "(in percent)" label in C2
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "(in percent)"

# ---------------------------------------------------------------------
# 2. Revise the existing 2022 share value (S5: 30 -> 29.9)
# ---------------------------------------------------------------------
$ws.Range("S5").Value = 29.9

# ---------------------------------------------------------------------
# 3. Re-width the year columns (D:T) to a uniform 9 characters, matching
#    the new single-width column band used for 2007-2023
# ---------------------------------------------------------------------
$ws.Range("D1:T1").EntireColumn.ColumnWidth = 8.1

# ---------------------------------------------------------------------
# 4. Add the new 2023 column (T) with header/year row, share row and
#    production row, copying the look of the neighbouring 2022 column (S)
# ---------------------------------------------------------------------

# -- T4: year header (matches style of S4) --
$t4 = $ws.Range("T4")
$t4.Font.Name = "Times New Roman"
$t4.Font.Size = 9
$t4.Font.Bold = $true
$t4.Borders.Item(8).LineStyle = 1
$t4.Borders.Item(8).Weight = -4138
$t4.Borders.Item(9).LineStyle = 1
$t4.Borders.Item(9).Weight = -4138
$t4.HorizontalAlignment = -4152
$t4.VerticalAlignment = -4108
$t4.WrapText = $true
$t4.Value = 2023

# -- T5: renewable energy share for 2023 (matches style of S5) --
$t5 = $ws.Range("T5")
$t5.Font.Name = "Times New Roman"
$t5.Font.Size = 9
$t5.Font.Bold = $false
$t5.HorizontalAlignment = -4152
$t5.VerticalAlignment = -4108
$t5.WrapText = $true
$t5.Value = 29.5

# -- T6: hydropower production for 2023 (matches style of S6) --
$t6 = $ws.Range("T6")
$t6.Font.Name = "Times New Roman"
$t6.Font.Size = 9
$t6.Font.Bold = $false
$t6.Borders.Item(9).LineStyle = 1
$t6.Borders.Item(9).Weight = -4138
$t6.HorizontalAlignment = -4152
$t6.VerticalAlignment = -4108
$t6.WrapText = $true
$t6.Value = 12030.6

# ---------------------------------------------------------------------
# 5. Reset the saved selection back to the default top-left cell
# ---------------------------------------------------------------------
$ws.Range("A1").Activate()
